$wb = $excel.ActiveWorkbook

# --- occurence.txt (sheet1): fix header typo, merge/rename event-date columns ---
$ws1 = $wb.Worksheets.Item("occurence.txt")

# Column I was "eventDateBeginning" (value 1964-07-28) and column J was
# "eventDateEnding" (empty). Rename I's header to "eventDate" and drop the
# now-redundant "eventDateEnding" column (J) entirely, shifting everything
# after it one column to the left.
$ws1.Range("I1").Value = "eventDate"
$ws1.Columns.Item(10).Delete()

# Fix the "occurenceID" typo -> "occurrenceID" (column C header)
$ws1.Range("C1").Value = "occurrenceID"

# The previously-active sheet (identification.txt) becomes inactive and
# occurence.txt becomes the selected / active tab with C1 selected.
$ws1.Activate()
$ws1.Range("C1").Select() | Out-Null
